$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The quarterly table (columns D:M, rows 8:58) rolls forward by one
# quarter: the oldest quarter (old column D, "1399/06") drops off and
# a new quarter ("1401/12") is appended as the new column M.
# ------------------------------------------------------------------

# 1) Drop the oldest quarter - this shifts D:M left by one column and
#    naturally carries over all the correct values/styles/widths for
#    the columns that remain (old E->D, old F->E, ... old M->L).
$ws.Columns("D").Delete()

# 2) Clone column L's formatting into the new column M so the new
#    column matches the existing style pattern exactly.
$ws.Range("L1:L59").Copy($ws.Range("M1:M59"))

# 3) Fill in the new quarter's header + values in column M.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# "M9" needs to hold literal text ("1402-02-28") rather than be
# auto-parsed as a date, so stage it on a scratch cell formatted as
# text, then paste just the value across (keeps M9's existing style).
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "1402-02-28"
$ws.Range("ZZ1").Copy()
$ws.Range("M9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("M12").Value = 1152716
$ws.Range("M14").Value = 7896235
$ws.Range("M15").Value = 3989200
$ws.Range("M16").Value = 1006201
$ws.Range("M18").Value = 14044352
$ws.Range("M19").Value = 12083
$ws.Range("M22").Value = 1611304
$ws.Range("M23").Value = 22957
$ws.Range("M26").Value = 1646555
$ws.Range("M27").Value = 15690907
$ws.Range("M29").Value = 1135654
$ws.Range("M31").Value = 86063
$ws.Range("M32").Value = 810071
$ws.Range("M33").Value = 5020189
$ws.Range("M34").Value = 2403520
$ws.Range("M37").Value = 9455497
$ws.Range("M41").Value = 445277
$ws.Range("M42").Value = 445277
$ws.Range("M43").Value = 9900774
$ws.Range("M56").Value = 5128350
$ws.Range("M57").Value = 5790133
$ws.Range("M58").Value = 15690907

# 4) Column M should render at the same width as the other
#    "publish-date row" columns (E, I, M = 31 chars wide).
$ws.Columns("M").ColumnWidth = 30.16

# ------------------------------------------------------------------
# One quarter's publish date was revised after the fact: what used to
# read "1401-11-01 (7)" (now sitting in column I after the shift)
# should read "1402-02-28 (8)".
# ------------------------------------------------------------------
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "1402-02-28 (8)"
$ws.Range("ZZ1").Copy()
$ws.Range("I9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
